$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21 (shifts existing rows 21-28 down to 22-29),
# inheriting the formatting of the surrounding rows.
$ws.Rows.Item(21).Insert()

# Fill in the new test case row (T10D). Set the Job Key (B) before the
# Test Name (A) so new shared-string entries are appended in the same
# order as the target workbook (T10D first, then the long test name).
$ws.Range("B21").Value = "T10D"
$ws.Range("A21").Value = "010-runtime-clustering - runtime-sharedfs-all-staging-nogridstart"
$ws.Range("C21").Value = "SharedFS"
$ws.Range("D21").Value = "Horizontal - Runtime"
$ws.Range("E21").Value = "Enabled"
$ws.Range("F21").Value = "Staged"
$ws.Range("G21").Value = "Disabled"

# The "blackdiamond - sharedfs-worker-staging-cleanup" row (now row 27)
# no longer lists a Worker/Gridstart value.
$ws.Range("G27").Clear()

# The "blackdiamond - horizontal" row (now row 28) gains a Gridstart value.
$ws.Range("G28").Value = "Disabled"

# The "blackdiamond - sharedfs-worker-staging-cleanup" duplicate row
# (now row 29) gains an Executable/Gridstart value.
$ws.Range("F29").Value = "Staged"

# Widen column A to fit the new, longer test name.
$ws.Columns.Item(1).ColumnWidth = 59.5

# Update the selection to match the edited range.
$ws.Range("A5:K29").Select()
